$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 <= old row 8
$ws.Range("A6").Value = 112126260
$ws.Range("B6").Value = 77515
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = "Garnlav"
$ws.Range("G6").Value = "Alectoria sarmentosa"
$ws.Range("H6").Value = "(Ach.) Ach."
$ws.Range("P6").Value = "Svarvarmyran (Svarvarmyran), Ly lm"
$ws.Range("Q6").Value = 690607.2411511086
$ws.Range("R6").Value = 7125723.544707977
$ws.Range("Z6").Value = "13:28"
$ws.Range("AB6").Value = "13:28"

# Row 7 <= old row 6
$ws.Range("A7").Value = 112129324
$ws.Range("B7").Value = 90666
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 4364
$ws.Range("F7").Value = "Dropptaggsvamp"
$ws.Range("G7").Value = "Hydnellum ferrugineum"
$ws.Range("H7").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q7").Value = 690443.61286689
$ws.Range("R7").Value = 7126178.415219921
$ws.Range("Z7").Value = "15:32"
$ws.Range("AB7").Value = "15:32"

# Row 8 <= old row 7
$ws.Range("A8").Value = 112129442
$ws.Range("B8").Value = 90710
$ws.Range("E8").Value = 5449
$ws.Range("F8").Value = "Svart taggsvamp"
$ws.Range("G8").Value = "Phellodon niger"
$ws.Range("H8").Value = "(Fr.:Fr.) P.Karst."
$ws.Range("P8").Value = "Godmyr (Godmyr), Ly lm"
$ws.Range("Q8").Value = 690472.3748497693
$ws.Range("R8").Value = 7126172.770408084
$ws.Range("Z8").Value = "15:36"
$ws.Range("AB8").Value = "15:36"

# Row 22 <= old row 25
$ws.Range("A22").Value = 112128981
$ws.Range("B22").Value = 78578
$ws.Range("E22").Value = 6458
$ws.Range("F22").Value = "Lunglav"
$ws.Range("G22").Value = "Lobaria pulmonaria"
$ws.Range("H22").Value = "(L.) Hoffm."
$ws.Range("P22").Value = "Godmyr (Godmyr), Ly lm"
$ws.Range("Q22").Value = 690353.6739480412
$ws.Range("R22").Value = 7126318.234514099
$ws.Range("Z22").Value = "15:13"
$ws.Range("AB22").Value = "15:13"

# Row 23 <= old row 22
$ws.Range("A23").Value = 112125962
$ws.Range("B23").Value = 90710
$ws.Range("D23").Value = "NT"
$ws.Range("E23").Value = 5449
$ws.Range("F23").Value = "Svart taggsvamp"
$ws.Range("G23").Value = "Phellodon niger"
$ws.Range("H23").Value = "(Fr.:Fr.) P.Karst."
$ws.Range("Q23").Value = 690606.1334164523
$ws.Range("R23").Value = 7125734.392117385
$ws.Range("Z23").Value = "13:22"
$ws.Range("AB23").Value = "13:22"

# Row 24 <= old row 23
$ws.Range("A24").Value = 112126910
$ws.Range("B24").Value = 90666
$ws.Range("D24").Value = "LC"
$ws.Range("E24").Value = 4364
$ws.Range("F24").Value = "Dropptaggsvamp"
$ws.Range("G24").Value = "Hydnellum ferrugineum"
$ws.Range("H24").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("P24").Value = "Svarvarmyran (Svarvarmyran), Ly lm"
$ws.Range("Q24").Value = 690565.1288052741
$ws.Range("R24").Value = 7125648.02276709
$ws.Range("Z24").Value = "13:54"
$ws.Range("AB24").Value = "13:54"

# Row 25 <= old row 24
$ws.Range("A25").Value = 112129248
$ws.Range("B25").Value = 90710
$ws.Range("E25").Value = 5449
$ws.Range("F25").Value = "Svart taggsvamp"
$ws.Range("G25").Value = "Phellodon niger"
$ws.Range("H25").Value = "(Fr.:Fr.) P.Karst."
$ws.Range("Q25").Value = 690368.2987376999
$ws.Range("R25").Value = 7126265.441720054
$ws.Range("Z25").Value = "15:28"
$ws.Range("AB25").Value = "15:28"

